$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new year column K, copying the formatting from column J (the last existing data column)
$ws.Range("J3:J6").Copy() | Out-Null
$ws.Range("K3").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Fill in the new values for 2023
$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 913.7
$ws.Range("K5").Value = 507.3
$ws.Range("K6").Value = 1068.5

$excel.CutCopyMode = 0
